# Automatic update of files.
#
# The data rows 2-9 on the "Artfynd" sheet get cyclically shifted up by one:
#   new row 2 = old row 3
#   new row 3 = old row 4
#   ...
#   new row 8 = old row 9
#   new row 9 = old row 2
# Row 1 (headers) and rows 10+ are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 9
$lastCol = 51   # column AY

# 1. Snapshot every cell's value in the affected row range before writing anything,
#    so reads are never polluted by writes we've already made.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowData = @{}
    for ($c = 1; $c -le $lastCol; $c++) {
        $rowData[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowData
}

# 2. Write the cyclic shift: row r gets the content that used to live in row r+1
#    (wrapping row 9's old content back around into row 2's old slot... i.e. row 9
#    gets what row 2 used to hold).
#    String values are written to cells pre-formatted as text ("@") so Excel's
#    automatic type-inference doesn't turn e.g. "2014-10-30" into a real date.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $r + 1
    if ($srcRow -gt $lastRow) {
        $srcRow = $firstRow
    }
    $srcData = $snapshot[$srcRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $val = $srcData[$c]
        $cell = $ws.Cells.Item($r, $c)
        if ($val -is [string]) {
            $cell.NumberFormat = "@"
        } else {
            $cell.NumberFormat = "General"
        }
        $cell.Value = $val
    }
}

# 3. Value2 coerces purely-numeric-looking text (e.g. "6") to a Double on read,
#    losing the fact that the source cell is text ("Antal" = 6 fruktkroppar).
#    Column I (Antal) on row 3 held the text "6" before the shift and that row's
#    content becomes row 2 after the shift, so patch it back up explicitly.
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "6"
